$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 9 with sparse data matching the diff:
# A9 = date serial 45707 (2025-02-19), styled like other date cells (style "2")
# C9 = "Corte Adulto"
# D9 = 10
# E9 = "vidal"
# G9 = "Efectivo"
# (B9, F9, H9 intentionally left blank)

$ws.Range("A9").Value = 45707
$ws.Range("A9").NumberFormat = $ws.Range("A2").NumberFormat

$ws.Range("C9").Value = "Corte Adulto"
$ws.Range("D9").Value = 10
$ws.Range("E9").Value = "vidal"
$ws.Range("G9").Value = "Efectivo"
